# Auto-generated edit script: updates cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRange, [string]$value) {
    # Force text storage for values that would otherwise be auto-parsed as a number
    # (mirrors a user typing an apostrophe prefix in Excel) so price strings like
    # '1.00' or '303.05' remain text cells, matching the sheet's existing convention.
    if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cellRange.Value = "'" + $value
    } else {
        $cellRange.Value = $value
    }
}

Set-CellText $ws.Range('D2') '42.702.00'
Set-CellText $ws.Range('E2') '  -1.20%  '
Set-CellText $ws.Range('D3') '2.308.08'
Set-CellText $ws.Range('E3') '  +0.24%  '
Set-CellText $ws.Range('D4') '1.00'
Set-CellText $ws.Range('E4') '  +0.04%  '
Set-CellText $ws.Range('D5') '303.05'
Set-CellText $ws.Range('E5') '  -1.85%  '
Set-CellText $ws.Range('D6') '99.76'
Set-CellText $ws.Range('E6') '  -4.18%  '
Set-CellText $ws.Range('D7') '0.506'
Set-CellText $ws.Range('E7') '  -3.50%  '
Set-CellText $ws.Range('D8') '1.00'
Set-CellText $ws.Range('E8') '  +0.13%  '
Set-CellText $ws.Range('D9') '0.503'
Set-CellText $ws.Range('E9') '  -3.22%  '
Set-CellText $ws.Range('D10') '34.77'
Set-CellText $ws.Range('E10') '  -2.94%  '
Set-CellText $ws.Range('D11') '0.0792'
Set-CellText $ws.Range('E11') '  -2.05%  '
Set-CellText $ws.Range('E12') '  +0.67%  '
Set-CellText $ws.Range('D13') '6.73'
Set-CellText $ws.Range('E13') '  -3.27%  '
Set-CellText $ws.Range('D14') '2.664.68'
Set-CellText $ws.Range('E14') '  +0.17%  '
Set-CellText $ws.Range('D15') '15.72'
Set-CellText $ws.Range('E15') '  +3.98%  '
Set-CellText $ws.Range('D16') '2.337.13'
Set-CellText $ws.Range('E16') '  +1.31%  '
Set-CellText $ws.Range('E17') '  +0.57%  '
Set-CellText $ws.Range('D18') '42.658.53'
Set-CellText $ws.Range('E18') '  -1.21%  '
Set-CellText $ws.Range('D19') '0.0₃0907'
Set-CellText $ws.Range('E19') '  -1.57%  '
Set-CellText $ws.Range('B20') 'InternetComputer(DFINITY)'
Set-CellText $ws.Range('C20') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws.Range('D20') '11.53'
Set-CellText $ws.Range('E20') '  -3.68%  '
Set-CellText $ws.Range('B21') 'Uniswap'
Set-CellText $ws.Range('C21') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText $ws.Range('D21') '6.07'
Set-CellText $ws.Range('E21') '  -1.41%  '
Set-CellText $ws.Range('D22') '67.99'
Set-CellText $ws.Range('E22') '  +0.31%  '
Set-CellText $ws.Range('D23') '235.30'
Set-CellText $ws.Range('E23') '  -2.04%  '
Set-CellText $ws.Range('D24') '1.96'
Set-CellText $ws.Range('E24') '  -2.00%  '
Set-CellText $ws.Range('D25') '2.53'
Set-CellText $ws.Range('E25') '  -2.62%  '
Set-CellText $ws.Range('E26') '  -0.12%  '
Set-CellText $ws.Range('D27') '24.95'
Set-CellText $ws.Range('E27') '  +1.18%  '
Set-CellText $ws.Range('E28') '  +3.15%  '
Set-CellText $ws.Range('D29') '34.67'
Set-CellText $ws.Range('E29') '  -4.18%  '
Set-CellText $ws.Range('D30') '164.23'
Set-CellText $ws.Range('E30') '  +1.72%  '
Set-CellText $ws.Range('D31') '9.15'
Set-CellText $ws.Range('E31') '  -4.23%  '
Set-CellText $ws.Range('D32') '1.00'
Set-CellText $ws.Range('E32') '  +0.08%  '
Set-CellText $ws.Range('D33') '5.02'
Set-CellText $ws.Range('E33') '  -4.25%  '
Set-CellText $ws.Range('E34') '  -4.87%  '
Set-CellText $ws.Range('E35') '  +0.15%  '
Set-CellText $ws.Range('D36') '16.82'
Set-CellText $ws.Range('E36') '  -7.66%  '
Set-CellText $ws.Range('D37') '0.0701'
Set-CellText $ws.Range('E37') '  -4.60%  '
Set-CellText $ws.Range('D38') '2.89'
Set-CellText $ws.Range('E38') '  -3.55%  '
Set-CellText $ws.Range('D39') '1.80'
Set-CellText $ws.Range('E39') '  -2.97%  '
Set-CellText $ws.Range('E40') '  -5.82%  '
Set-CellText $ws.Range('E41') '  -3.48%  '
Set-CellText $ws.Range('D42') '2.49'
Set-CellText $ws.Range('E42') '  -1.63%  '
Set-CellText $ws.Range('D43') '1.971.64'
Set-CellText $ws.Range('E43') '  +0.48%  '
Set-CellText $ws.Range('D44') '0.0280'
Set-CellText $ws.Range('E44') '  -2.94%  '
Set-CellText $ws.Range('D45') '18.49'
Set-CellText $ws.Range('E45') '  -1.54%  '
Set-CellText $ws.Range('D46') '10.24'
Set-CellText $ws.Range('E46') '  +0.57%  '
Set-CellText $ws.Range('D47') '2.89'
Set-CellText $ws.Range('E47') '  -5.87%  '
Set-CellText $ws.Range('D48') '55.71'
Set-CellText $ws.Range('E48') '  -2.38%  '
Set-CellText $ws.Range('B49') 'RocketPoolETH'
Set-CellText $ws.Range('C49') 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-CellText $ws.Range('D49') '2.532.95'
Set-CellText $ws.Range('E49') '  +0.16%  '
Set-CellText $ws.Range('B50') 'HuobiToken'
Set-CellText $ws.Range('C50') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText $ws.Range('D50') '2.84'
Set-CellText $ws.Range('E50') '  -3.08%  '
Set-CellText $ws.Range('E51') '  +0.53%  '
